# Applies the "Vaccine" price-list text fixes described in the commit:
#   - normalize several vaccine-name strings (typos, missing slashes/spaces)
#   - split the merged "packaging" text for ENGERIX B (rows 27-30) into the
#     four distinct packaging descriptions it always should have had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Vaccine-name (column A) text corrections -----------------------------

$ws.Range("A2").Value  = "DTaP/"
$ws.Range("A3").Value  = "DTaP/"
$ws.Range("A4").Value  = "DTaP/"
$ws.Range("A5").Value  = "DTaP/"

$ws.Range("A8").Value  = "DTaP-Hib "

$ws.Range("A10").Value = "Hepatitis B-Hib"

$ws.Range("A21").Value = "Hepatitis A-Hepatitis B Adult"
$ws.Range("A22").Value = "Hepatitis A-Hepatitis B Adult"
$ws.Range("A23").Value = "Hepatitis A-Hepatitis B Adult"

$ws.Range("A24").Value = "Hepatitis A-Hepatitis B 18 only"
$ws.Range("A25").Value = "Hepatitis A-Hepatitis B 18 only"
$ws.Range("A26").Value = "Hepatitis A-Hepatitis B 18 only"

$ws.Range("A27").Value = "Hepatitis B Pediatric/Adolescent"
$ws.Range("A28").Value = "Hepatitis B Pediatric/Adolescent"
$ws.Range("A29").Value = "Hepatitis B Pediatric/Adolescent"
$ws.Range("A30").Value = "Hepatitis B Pediatric/Adolescent"
$ws.Range("A31").Value = "Hepatitis B Pediatric/Adolescent"

$ws.Range("A43").Value = "Influenza (Live, Intranasal)"

$ws.Range("A44").Value = "MMR/"

$ws.Range("A45").Value = "Pneumococcal 7-valent (Pediatric)"

$ws.Range("A47").Value = "Tetanus  Diphtheria Toxoids"
$ws.Range("A48").Value = "Tetanus  Diphtheria Toxoids"

# --- Packaging (column D) split for ENGERIX B, rows 27-30 -----------------
# Previously all four rows shared one run-on string; split into the four
# distinct packaging sizes.

$ws.Range("D27").Value = "1 dose vials "
$ws.Range("D28").Value = "10 pack - 1 dose vials "
$ws.Range("D29").Value = "5 pack - 1 dose T-L syringes, No Needle "
$ws.Range("D30").Value = "25 pack - 1 dose T-L syringes, No Needle "
